$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.281.50"
$ws.Range("E2").Value = "  +3.90%  "
$ws.Range("D3").Value = "1.728.08"
$ws.Range("E3").Value = "  +2.27%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.523"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.92"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.78%  "
$ws.Range("E9").Value = "  +2.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0638"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.54%  "
$ws.Range("E11").Value = "  +0.53%  "
$ws.Range("D12").Value = "1.972.45"
$ws.Range("E12").Value = "  +2.30%  "
$ws.Range("D13").Value = "1.727.93"
$ws.Range("E13").Value = "  +2.21%  "
$ws.Range("E14").Value = "  +0.89%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.565"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.74%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.65"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("D17").Value = "28.259.84"
$ws.Range("E17").Value = "  +3.80%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "244.79"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.10%  "
$ws.Range("E19").Value = "  +1.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.93"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.41%  "
$ws.Range("E21").Value = "  -0.19%  "
$ws.Range("E22").Value = "  +1.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.67"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.42%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.08"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.38"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.03%  "
$ws.Range("E26").Value = "  +2.51%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.65"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.05%  "
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("E29").Value = "  -0.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0517"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.19%  "
$ws.Range("E31").Value = "  +2.38%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.43"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.75%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.26"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.82%  "
$ws.Range("D34").Value = "1.481.80"
$ws.Range("E34").Value = "  -4.46%  "
$ws.Range("E35").Value = "  -2.48%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.980"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.71%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.605"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.40"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.07%  "
$ws.Range("E39").Value = "  +0.88%  "
$ws.Range("E40").Value = "  +0.27%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "69.74"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.79%  "
$ws.Range("E42").Value = "  -0.23%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.64"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.61%  "
$ws.Range("D44").Value = "1.876.24"
$ws.Range("E44").Value = "  +2.07%  "
$ws.Range("E45").Value = "  +1.40%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.801"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.59%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.72"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.26%  "
$ws.Range("D48").Value = "0.0₆0113"
$ws.Range("E48").Value = "  +2.97%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "90.19"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.16%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.15"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.90%  "
$ws.Range("E51").Value = "  -0.83%  "
